$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("E1").Value = "Identifier"

$ws.Range("E2").Formula = "=CONCATENATE(C2,"" "",D2)"
$ws.Range("E3").Formula = "=CONCATENATE(C3,"" "",D3)"
$ws.Range("E4").Formula = "=CONCATENATE(C4,"" "",D4)"
$ws.Range("E5").Formula = "=CONCATENATE(C5,"" "",D5)"
$ws.Range("E6").Formula = "=CONCATENATE(C6,"" "",D6)"
$ws.Range("E7").Formula = "=CONCATENATE(C7,"" "",D7)"

$ws.Columns.Item(3).ColumnWidth = 13.2422
$ws.Columns.Item(5).ColumnWidth = 13.5
